# Update "想去人数" (number of people interested) figures for a handful of
# events across the "展览" (Exhibitions), "演出" (Shows) and "全部类型"
# (All types) worksheets, matching the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 307
$wsExpo.Range("F4").Value = 8103
$wsExpo.Range("F5").Value = 5909
$wsExpo.Range("F10").Value = 294
$wsExpo.Range("F11").Value = 444

# --- 演出 (Shows) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 7

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 307
$wsAll.Range("F4").Value = 8103
$wsAll.Range("F5").Value = 5909
$wsAll.Range("F10").Value = 294
$wsAll.Range("F14").Value = 7
$wsAll.Range("F15").Value = 444
